$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.249.77"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "3.519.12"
$ws.Range("E3").Value = "  -0.33%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.52"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.06"
$ws.Range("E6").Value = "  +0.49%  "

# Row 7
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.514.84"
$ws.Range("E7").Value = "  -0.40%  "

# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  -1.91%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.195"
$ws.Range("E10").Value = "  -1.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.15"
$ws.Range("E11").Value = "  +5.43%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.587"
$ws.Range("E12").Value = "  +0.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.15"
$ws.Range("E13").Value = "  -1.85%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000276"
$ws.Range("E14").Value = "  -0.84%  "

# Row 15
$ws.Range("D15").Value = "4.081.33"
$ws.Range("E15").Value = "  -0.31%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  +0.36%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "611.49"
$ws.Range("E17").Value = "  -2.27%  "

# Row 18
$ws.Range("D18").Value = "3.518.58"
$ws.Range("E18").Value = "  -0.34%  "

# Row 19
$ws.Range("D19").Value = "70.272.11"
$ws.Range("E19").Value = "  -0.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.122"
$ws.Range("E20").Value = "  +2.21%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.73"
$ws.Range("E21").Value = "  +1.96%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.880"
$ws.Range("E22").Value = "  -0.49%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.12"
$ws.Range("E23").Value = "  -4.99%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "99.16"
$ws.Range("E24").Value = "  +3.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.55"
$ws.Range("E25").Value = "  -2.13%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.74"
$ws.Range("E26").Value = "  -2.36%  "

# Row 27
$ws.Range("E27").Value = "  -0.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.58"
$ws.Range("E28").Value = "  -1.22%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.86"
$ws.Range("E29").Value = "  +1.53%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.09"
$ws.Range("E30").Value = "  -0.89%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.12"
$ws.Range("E31").Value = "  -4.17%  "

# Row 32
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.99"
$ws.Range("E32").Value = "  -2.46%  "

# Row 33
$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.31"
$ws.Range("E33").Value = "  -0.60%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.80"
$ws.Range("E34").Value = "  -1.88%  "

# Row 35
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "618.07"
$ws.Range("E35").Value = "  +7.78%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0506"
$ws.Range("E36").Value = "  +10.13%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.100"
$ws.Range("E37").Value = "  -1.98%  "

# Row 38
$ws.Range("B38").Value = "Cosmos"
$ws.Range("C38").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.77"
$ws.Range("E38").Value = "  -0.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.46"
$ws.Range("E39").Value = "  -3.23%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.147"
$ws.Range("E40").Value = "  +2.40%  "

# Row 41
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.77"
$ws.Range("E41").Value = "  -1.34%  "

# Row 42
$ws.Range("E42").Value = "  -0.10%  "

# Row 43
$ws.Range("D43").Value = "3.377.64"
$ws.Range("E43").Value = "  +0.88%  "

# Row 44
$ws.Range("D44").Value = "0.0₃0742"
$ws.Range("E44").Value = "  +4.86%  "

# Row 45
$ws.Range("E45").Value = "  -5.04%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.05"
$ws.Range("E46").Value = "  -3.27%  "

# Row 47
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.88"
$ws.Range("E47").Value = "  -6.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.56"
$ws.Range("E48").Value = "  -3.87%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.131"
$ws.Range("E49").Value = "  +0.61%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.54"
$ws.Range("E50").Value = "  -0.68%  "
